# Update New Orleans xlsx files:
#  1. Insert a new "State" column into hotel_info (after Hotel_Name, before City)
#     with value "Louisiana" for the existing hotel row.
#  2. Reorder the worksheets so "review_info" comes before "hotel_info".

$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")

# Insert a new column C (State), shifting City/Zip/etc one column to the right.
$hotelSheet.Columns("C:C").Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# Move review_info to be the first sheet (tab order), hotel_info becomes second.
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($wb.Worksheets.Item(1))
